$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.485.38"
$ws.Range("E2").Value = "  +4.33%  "

# Row 3
$ws.Range("D3").Value = "3.608.35"
$ws.Range("E3").Value = "  +4.12%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "'629.46"
$ws.Range("E5").Value = "  +4.43%  "

# Row 6
$ws.Range("D6").Value = "'158.80"
$ws.Range("E6").Value = "  +7.06%  "

# Row 7
$ws.Range("D7").Value = "3.605.53"
$ws.Range("E7").Value = "  +3.96%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  +3.32%  "

# Row 10
$ws.Range("E10").Value = "  +8.83%  "

# Row 11
$ws.Range("D11").Value = "'7.41"
$ws.Range("E11").Value = "  +7.57%  "

# Row 12
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  +4.67%  "

# Row 13
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +5.64%  "

# Row 14
$ws.Range("D14").Value = "'33.54"
$ws.Range("E14").Value = "  +7.09%  "

# Row 15
$ws.Range("D15").Value = "4.214.65"
$ws.Range("E15").Value = "  +3.72%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.610.56"
$ws.Range("E16").Value = "  +3.83%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.522.06"
$ws.Range("E17").Value = "  +4.32%  "

# Row 18
$ws.Range("E18").Value = "  +0.69%  "

# Row 19
$ws.Range("E19").Value = "  +5.52%  "

# Row 20
$ws.Range("D20").Value = "'16.10"
$ws.Range("E20").Value = "  +7.44%  "

# Row 21
$ws.Range("D21").Value = "'10.29"
$ws.Range("E21").Value = "  +14.61%  "

# Row 22
$ws.Range("D22").Value = "'463.24"
$ws.Range("E22").Value = "  +4.67%  "

# Row 23
$ws.Range("D23").Value = "'0.647"
$ws.Range("E23").Value = "  +4.16%  "

# Row 24
$ws.Range("D24").Value = "'78.63"
$ws.Range("E24").Value = "  +1.94%  "

# Row 25
$ws.Range("E25").Value = "  +9.29%  "

# Row 26
$ws.Range("D26").Value = "'10.72"
$ws.Range("E26").Value = "  +6.97%  "

# Row 27
$ws.Range("D27").Value = "3.748.72"
$ws.Range("E27").Value = "  +3.82%  "

# Row 28
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  +12.79%  "

# Row 30
$ws.Range("E30").Value = "  +5.24%  "

# Row 31
$ws.Range("D31").Value = "'1.71"
$ws.Range("E31").Value = "  +9.58%  "

# Row 32
$ws.Range("D32").Value = "'0.175"
$ws.Range("E32").Value = "  +9.63%  "

# Row 33
$ws.Range("D33").Value = "'6.56"
$ws.Range("E33").Value = "  +7.35%  "

# Row 34
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
$ws.Range("D35").Value = "'26.57"
$ws.Range("E35").Value = "  +4.28%  "

# Row 36
$ws.Range("E36").Value = "  +5.32%  "

# Row 37
$ws.Range("D37").Value = "3.601.79"
$ws.Range("E37").Value = "  +4.08%  "

# Row 38
$ws.Range("D38").Value = "'8.52"
$ws.Range("E38").Value = "  +7.77%  "

# Row 39
$ws.Range("D39").Value = "'2.42"
$ws.Range("E39").Value = "  +11.93%  "

# Row 40
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0927"
$ws.Range("E41").Value = "  +7.64%  "

# Row 42
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'179.71"
$ws.Range("E42").Value = "  +3.74%  "

# Row 43
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.26%  "

# Row 44
$ws.Range("E44").Value = "  +3.75%  "

# Row 45
$ws.Range("D45").Value = "'31.80"
$ws.Range("E45").Value = "  +22.28%  "

# Row 46
$ws.Range("E46").Value = "  +3.93%  "

# Row 47
$ws.Range("E47").Value = "  +13.70%  "

# Row 48
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'45.91"
$ws.Range("E48").Value = "  +1.57%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.75"
$ws.Range("E49").Value = "  +11.02%  "

# Row 50
$ws.Range("D50").Value = "'7.83"
$ws.Range("E50").Value = "  +4.03%  "

# Row 51
$ws.Range("E51").Value = "  +8.99%  "
